$wb = $excel.ActiveWorkbook

# --- Sheet "All Orders": update row 18 ---
$ws1 = $wb.Worksheets.Item("All Orders")
$ws1.Range("H18").Value = "CANCELLED"
$ws1.Range("M18").Value = "test order"

# --- Sheet "Daily Summary": update row 4 totals ---
$ws2 = $wb.Worksheets.Item("Daily Summary")
$ws2.Range("D4").Value = 5
$ws2.Range("E4").Value = 155
$ws2.Range("G4").Value = 155
